# Updates "horarios" workbook (3 sheets) with the newer scrape snapshot
# (Linea 141 - 1158, timestamp 11:54:47), per commit "Horarios actualizados
# Linea 141 - 1158".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 11:54:47"
$ws1.Range("A3").Value = "Total filas: 162"

# Small in-place corrections scattered through the already-existing rows
# (the underlying scrape re-ordered a couple of near-simultaneous arrivals).
$ws1.Range("C22").Value = "215C_EL PATO"
$ws1.Range("C23").Value = "14_ABASTO"

$ws1.Range("A37").Value = "07:28:23"
$ws1.Range("C37").Value = "16_SANTA ANA"
$ws1.Range("D37").Value = 3

$ws1.Range("A38").Value = "06:55:48"
$ws1.Range("C38").Value = "11_ETCHEVERRY"
$ws1.Range("D38").Value = 36

$ws1.Range("A62").Value = "08:41:16"
$ws1.Range("C62").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D62").Value = 2

$ws1.Range("A63").Value = "08:04:39"
$ws1.Range("C63").Value = "14_ABASTO"
$ws1.Range("D63").Value = 39

$ws1.Range("C68").Value = "11_ETCHEVERRY"
$ws1.Range("C70").Value = "215A_EL PATO"

$ws1.Range("C76").Value = "16_SANTA ANA"
$ws1.Range("C77").Value = "16_P MOR-SANTA ANA"

$ws1.Range("C82").Value = "16_SANTA ANA"
$ws1.Range("C83").Value = "17_ROMERO"

# Rows 137-167: newest batch of scraped arrivals (new run timestamp
# 11:54:47) replaces/extends what used to be rows 137-150.
$rows1 = @(
    @(137,"11:54:47","11:54","16_SANTA ANA",0,"LP1912"),
    @(138,"11:54:47","11:56","16_SANTA ANA",2,"LP1912"),
    @(139,"11:54:47","11:59","225_GOMEZ",5,"LP1912"),
    @(140,"11:54:47","12:02","84_COLONIA URQUIZA-ESC 49",8,"LP1912"),
    @(141,"11:07:42","12:06","14_ABASTO",59,"LP1912"),
    @(142,"11:54:47","12:06","23_HERNANDEZ",12,"LP1912"),
    @(143,"11:54:47","12:06","16_P MOR-SANTA ANA",12,"LP1912"),
    @(144,"10:20:05","12:07","16_P MOR-SANTA ANA",107,"LP1912"),
    @(145,"11:07:42","12:08","10_OLMOS",61,"LP1912"),
    @(146,"11:54:47","12:13","10_OLMOS",19,"LP1912"),
    @(147,"11:54:47","12:14","17_ROMERO",20,"LP1912"),
    @(148,"10:20:05","12:18","14_ABASTO",118,"LP1912"),
    @(149,"11:54:47","12:20","14_ABASTO",26,"LP1912"),
    @(150,"11:54:47","12:20","215A_EL PATO",26,"LP1912"),
    @(151,"11:54:47","12:21","26_HERNANDEZ",27,"LP1912"),
    @(152,"11:54:47","12:34","11_ETCHEVERRY",40,"LP1912"),
    @(153,"11:54:47","12:34","23_HERNANDEZ",40,"LP1912"),
    @(154,"11:54:47","12:36","27_EL RETIRO",42,"LP1912"),
    @(155,"11:54:47","12:38","17_179 Y 38",44,"LP1912"),
    @(156,"11:54:47","12:41","10_OLMOS",47,"LP1912"),
    @(157,"11:54:47","12:48","11_ETCHEVERRY",54,"LP1912"),
    @(158,"11:54:47","12:55","10_OLMOS",61,"LP1912"),
    @(159,"11:54:47","13:06","16_P MOR-SANTA ANA",72,"LP1912"),
    @(160,"11:54:47","13:13","215D_EL PATO",79,"LP1912"),
    @(161,"11:54:47","13:20","10_OLMOS",86,"LP1912"),
    @(162,"11:54:47","13:21","26_HERNANDEZ",87,"LP1912"),
    @(163,"11:54:47","13:26","15_ABASTO",92,"LP1912"),
    @(164,"11:54:47","13:26","14_ABASTO",92,"LP1912"),
    @(165,"11:54:47","13:36","15_ABASTO",102,"LP1912"),
    @(166,"11:54:47","13:46","17_ROMERO",112,"LP1912"),
    @(167,"11:54:47","13:50","215A_EL PATO",116,"LP1912")
)

foreach ($row in $rows1) {
    $r = $row[0]
    $ws1.Cells.Item($r,1).Value = $row[1]
    $ws1.Cells.Item($r,2).Value = $row[2]
    $ws1.Cells.Item($r,3).Value = $row[3]
    $ws1.Cells.Item($r,4).Value = $row[4]
    $ws1.Cells.Item($r,5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 11:54:47"
$ws2.Range("A3").Value = "Total filas: 18"

$rows2 = @(
    @(21,"11:54:47","12:20","215A_EL PATO",26,"LP1912"),
    @(22,"11:54:47","13:13","215D_EL PATO",79,"LP1912"),
    @(23,"11:54:47","13:50","215A_EL PATO",116,"LP1912")
)

foreach ($row in $rows2) {
    $r = $row[0]
    $ws2.Cells.Item($r,1).Value = $row[1]
    $ws2.Cells.Item($r,2).Value = $row[2]
    $ws2.Cells.Item($r,3).Value = $row[3]
    $ws2.Cells.Item($r,4).Value = $row[4]
    $ws2.Cells.Item($r,5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 11:54:47"
$ws3.Range("A3").Value = "Total filas: 25"

$rows3 = @(
    @(28,"11:54:47","12:04","215A_LA PLATA",10,"L6173"),
    @(29,"11:54:47","12:53","215C_LA PLATA",59,"L6203"),
    @(30,"11:54:47","13:30","215B_LP-P MOR-1 Y 57",96,"L6173")
)

foreach ($row in $rows3) {
    $r = $row[0]
    $ws3.Cells.Item($r,1).Value = $row[1]
    $ws3.Cells.Item($r,2).Value = $row[2]
    $ws3.Cells.Item($r,3).Value = $row[3]
    $ws3.Cells.Item($r,4).Value = $row[4]
    $ws3.Cells.Item($r,5).Value = $row[5]
}
